$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row (row 62) with the next day's data
$ws.Range("A62").Value = 46011
$ws.Range("B62").Value = 131
$ws.Range("C62").Value = 147
$ws.Range("D62").Value = 138

# Match the style of the preceding date cell (A61) for the new date cell (A62)
$ws.Range("A61").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$excel.CutCopyMode = $false
